$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 225, shifting existing rows 225:334 down to 226:335
$ws.Rows.Item(225).Insert()

# Populate the newly inserted row 225 with the new record's data
$ws.Range("A225").Value = 3
$ws.Range("B225").Value = "Femacal de La Calera"
$ws.Range("C225").Value = "Coquimbo"
$ws.Range("D225").Value = 44572
$ws.Range("E225").Value = 5
$ws.Range("F225").Value = 100112003
$ws.Range("G225").Value = "Ajo"
$ws.Range("H225").Value = "Chino"
$ws.Range("I225").Value = "Primera"
$ws.Range("J225").Value = 85
$ws.Range("K225").Value = 16000
$ws.Range("L225").Value = 16500
$ws.Range("M225").Value = 16265
$ws.Range("N225").Value = "$/caja 10 kilos"
$ws.Range("O225").Value = "Llay Llay"
$ws.Range("P225").Value = 1626
$ws.Range("Q225").Value = 10
$ws.Range("R225").Value = "Hortaliza"

# Copy the date-number style (s="2") used throughout column D onto the new cell
$ws.Range("D226").Copy()
$ws.Range("D225").PasteSpecial(-4122)
